$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 388.75
$ws.Range("I4").Value = 185
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 185
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -71
$ws.Range("N4").Value = -1228

$ws.Range("H5").Value = 111.15385
$ws.Range("I5").Value = 106.44444
$ws.Range("K5").Value = 106.44444
$ws.Range("M5").Value = 8.55556

$ws.Range("H8").Value = 2222402.2
$ws.Range("I8").Value = 2222402.2
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 6667206.600000001
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -6667067.600000001
$ws.Range("N8").ClearContents()

$ws.Range("H52").Value = 642.5
$ws.Range("I52").Value = 200
$ws.Range("J52").Value = 790
$ws.Range("K52").Value = 600
$ws.Range("L52").Value = 2370
$ws.Range("M52").Value = -440
$ws.Range("N52").Value = -2690

$ws.Range("H70").Value = 3764.5
$ws.Range("I70").Value = 1798
$ws.Range("K70").Value = 5394
$ws.Range("M70").Value = -5124

$ws.Range("H73").Value = 3764.5
$ws.Range("I73").Value = 1798
$ws.Range("K73").Value = 5394
$ws.Range("M73").Value = -4458

$ws.Range("H96").Value = 1165.6154
$ws.Range("I96").Value = 741.4286
$ws.Range("J96").Value = 1660.5
$ws.Range("K96").Value = 2224.2858
$ws.Range("L96").Value = 4981.5
$ws.Range("M96").Value = -851.2857999999997
$ws.Range("N96").Value = -7727.5

$ws.Range("H107").Value = 2142.7
$ws.Range("I107").Value = 2090.875
$ws.Range("J107").Value = 2350
$ws.Range("K107").Value = 2090.875
$ws.Range("L107").Value = 2350
$ws.Range("M107").Value = -170.875
$ws.Range("N107").Value = -6190

$ws.Range("H113").Value = 4919.5
$ws.Range("J113").Value = 4107.5
$ws.Range("L113").Value = 4107.5
$ws.Range("N113").Value = -10615.5

$ws.Range("H132").Value = 14824.25
$ws.Range("I132").Value = 13499.0625
$ws.Range("K132").Value = 40497.1875
$ws.Range("M132").Value = -37967.1875

$ws.Range("H137").Value = 1414.0834
$ws.Range("I137").Value = 1129.8334
$ws.Range("J137").Value = 1698.3334
$ws.Range("K137").Value = 3389.5002
$ws.Range("L137").Value = 5095.0002
$ws.Range("M137").Value = -839.5001999999999
$ws.Range("N137").Value = -10195.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1750.9286
$ws.Range("I45").Value = 1016.5455
$ws.Range("K45").Value = 1016.5455
$ws.Range("M45").Value = -639.5454999999999

$ws.Range("H132").Value = 4057.862
$ws.Range("I132").Value = 4042.75
$ws.Range("K132").Value = 12128.25
$ws.Range("M132").Value = -9598.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2005
$ws.Range("I105").Value = 2005
$ws.Range("K105").Value = 2005
$ws.Range("M105").Value = -258

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 3336.4375
$ws.Range("I7").Value = 5662.1113
$ws.Range("J7").Value = 346.2857
$ws.Range("K7").Value = 5662.1113
$ws.Range("L7").Value = 346.2857
$ws.Range("M7").Value = -5549.1113
$ws.Range("N7").Value = -572.2857

$ws.Range("H22").Value = 4250
$ws.Range("I22").Value = 4000
$ws.Range("K22").Value = 4000
$ws.Range("M22").Value = -3650

$ws.Range("H132").Value = 1996.75
$ws.Range("I132").Value = 1996.75
$ws.Range("K132").Value = 5990.25
$ws.Range("M132").Value = -3460.25

$ws.Range("H134").Value = 2156.5293
$ws.Range("I134").Value = 2156.5293
$ws.Range("K134").Value = 6469.5879
$ws.Range("M134").Value = -3934.5879

$ws.Range("H137").Value = 100000
$ws.Range("J137").Value = 100000
$ws.Range("L137").Value = 100000
$ws.Range("N137").Value = -110200

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 89.40000000000001
$ws.Range("I10").Value = 36.75
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 110.25
$ws.Range("L10").Value = 900
$ws.Range("M10").Value = 28.75
$ws.Range("N10").Value = -1178

$ws.Range("H17").Value = 9673.1
$ws.Range("I17").Value = 500.5
$ws.Range("J17").Value = 11966.25
$ws.Range("K17").Value = 1501.5
$ws.Range("L17").Value = 35898.75
$ws.Range("M17").Value = -1332.5
$ws.Range("N17").Value = -36236.75

$ws.Range("H38").Value = 348.25
$ws.Range("I38").Value = 345.85715
$ws.Range("K38").Value = 1037.57145
$ws.Range("M38").Value = -690.5714499999999

$ws.Range("H43").Value = 100
$ws.Range("I43").Value = 100
$ws.Range("K43").Value = 300
$ws.Range("M43").Value = -186

$ws.Range("H86").Value = 500
$ws.Range("I86").Value = 500
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -314
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 500
$ws.Range("I89").Value = 500
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 4500
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 1428
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4043.7778
$ws.Range("I70").Value = 3985.4285
$ws.Range("K70").Value = 3985.4285
$ws.Range("M70").Value = -3715.4285

$ws.Range("H73").Value = 4043.7778
$ws.Range("I73").Value = 3985.4285
$ws.Range("K73").Value = 3985.4285
$ws.Range("M73").Value = -3049.4285

$ws.Range("H122").Value = 3540.6667
$ws.Range("I122").Value = 2693.4
$ws.Range("J122").Value = 7777
$ws.Range("K122").Value = 8080.200000000001
$ws.Range("L122").Value = 23331
$ws.Range("M122").Value = -5630.200000000001
$ws.Range("N122").Value = -28231

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6712.857
$ws.Range("J46").Value = 7681.6665
$ws.Range("L46").Value = 7681.6665
$ws.Range("N46").Value = -8057.6665

$ws.Range("H68").Value = 4240.909
$ws.Range("J68").Value = 5642.857
$ws.Range("L68").Value = 5642.857
$ws.Range("N68").Value = -7140.857

$ws.Range("H71").Value = 4240.909
$ws.Range("J71").Value = 5642.857
$ws.Range("L71").Value = 28214.285
$ws.Range("N71").Value = -35702.285

$ws.Range("H93").Value = 991.6667
$ws.Range("J93").Value = 985
$ws.Range("L93").Value = 985
$ws.Range("N93").Value = -3481

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 16948.334
$ws.Range("I4").Value = 33816.668
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 33816.668
$ws.Range("L4").Value = 80
$ws.Range("M4").Value = -33703.668
$ws.Range("N4").Value = -306

$ws.Range("H107").Value = 1517.8889
$ws.Range("J107").Value = 616.6667
$ws.Range("L107").Value = 1850.0001
$ws.Range("N107").Value = -5690.0001

$ws.Range("H122").Value = 5964.143
$ws.Range("I122").Value = 3737.25
$ws.Range("J122").Value = 8933.333000000001
$ws.Range("K122").Value = 11211.75
$ws.Range("L122").Value = 26799.999
$ws.Range("M122").Value = -8761.75
$ws.Range("N122").Value = -31699.999

$ws.Range("H136").Value = 3963.5908
$ws.Range("I136").Value = 2680.818
$ws.Range("J136").Value = 5246.364
$ws.Range("K136").Value = 8042.454000000001
$ws.Range("L136").Value = 15739.092
$ws.Range("M136").Value = -5492.454000000001
$ws.Range("N136").Value = -20839.092
